# Auto-generated edits applying the diff to Mateus_Profits.xlsx (multi-sheet workbook)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 161.66667
$ws.Range("I9").Value = 161.66667
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 161.66667
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 7.333329999999989
$ws.Range("N9").ClearContents()

$ws.Range("H17").Value = 9093823
$ws.Range("J17").Value = 9093823
$ws.Range("L17").Value = 27281469
$ws.Range("N17").Value = -27281805

$ws.Range("H19").Value = 1747.1034
$ws.Range("I19").Value = 1202.8667
$ws.Range("K19").Value = 1202.8667
$ws.Range("M19").Value = -1027.8667

$ws.Range("H62").Value = 9411.444
$ws.Range("I62").Value = 12900.8
$ws.Range("K62").Value = 12900.8
$ws.Range("M62").Value = -12276.8

$ws.Range("H65").Value = 9411.444
$ws.Range("I65").Value = 12900.8
$ws.Range("K65").Value = 64504
$ws.Range("M65").Value = -61384

$ws.Range("H125").Value = 1846006.4
$ws.Range("I125").Value = 2584109
$ws.Range("J125").Value = 750
$ws.Range("K125").Value = 23256981
$ws.Range("L125").Value = 6750
$ws.Range("M125").Value = -23254521
$ws.Range("N125").Value = -11670

$ws.Range("H135").Value = 816.3214
$ws.Range("I135").Value = 841.5185
$ws.Range("K135").Value = 7573.6665
$ws.Range("M135").Value = -5038.6665

$ws.Range("H137").Value = 2570.1482
$ws.Range("I137").Value = 1991.7273
$ws.Range("J137").Value = 5115.2
$ws.Range("K137").Value = 5975.1819
$ws.Range("L137").Value = 15345.6
$ws.Range("M137").Value = -3425.1819
$ws.Range("N137").Value = -20445.6

$ws.Range("H138").Value = 4746.607
$ws.Range("I138").Value = 5696.3
$ws.Range("J138").Value = 4219
$ws.Range("K138").Value = 17088.9
$ws.Range("L138").Value = 12657
$ws.Range("M138").Value = -11948.9
$ws.Range("N138").Value = -22937

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7942.975
$ws.Range("I32").Value = 7942.975
$ws.Range("K32").Value = 7942.975
$ws.Range("M32").Value = -7655.975

$ws.Range("H42").Value = 4928
$ws.Range("I42").Value = 4928
$ws.Range("K42").Value = 4928
$ws.Range("M42").Value = -4442

$ws.Range("H132").Value = 1813.0358
$ws.Range("J132").Value = 1689
$ws.Range("L132").Value = 5067
$ws.Range("N132").Value = -10127

$ws.Range("H140").Value = 129999
$ws.Range("J140").Value = 129999
$ws.Range("L140").Value = 129999
$ws.Range("N140").Value = -140359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5051.0557
$ws.Range("I20").Value = 4977.304
$ws.Range("J20").Value = 5181.5386
$ws.Range("K20").Value = 4977.304
$ws.Range("L20").Value = 5181.5386
$ws.Range("M20").Value = -4730.304
$ws.Range("N20").Value = -5675.5386

$ws.Range("H99").Value = 1574.6923
$ws.Range("I99").Value = 1392.2632
$ws.Range("J99").Value = 2069.8572
$ws.Range("K99").Value = 1392.2632
$ws.Range("L99").Value = 2069.8572
$ws.Range("M99").Value = 105.7367999999999
$ws.Range("N99").Value = -5065.8572

$ws.Range("H105").Value = 3781.8572
$ws.Range("I105").Value = 3781.8572
$ws.Range("K105").Value = 3781.8572
$ws.Range("M105").Value = -2034.8572

$ws.Range("H134").Value = 4314.5454
$ws.Range("I134").Value = 3565.2068
$ws.Range("K134").Value = 10695.6204
$ws.Range("M134").Value = -8160.6204

$ws.Range("H140").Value = 79494.5
$ws.Range("J140").Value = 79494.5
$ws.Range("L140").Value = 79494.5
$ws.Range("N140").Value = -89854.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4872.4736
$ws.Range("I16").Value = 3680.0908
$ws.Range("K16").Value = 3680.0908
$ws.Range("M16").Value = -3393.0908

$ws.Range("H64").Value = 41956.168
$ws.Range("J64").Value = 41956.168
$ws.Range("L64").Value = 41956.168
$ws.Range("N64").Value = -42452.168

$ws.Range("H67").Value = 41956.168
$ws.Range("J67").Value = 41956.168
$ws.Range("L67").Value = 41956.168
$ws.Range("N67").Value = -43672.168

$ws.Range("H113").Value = 4872.4736
$ws.Range("I113").Value = 3680.0908
$ws.Range("K113").Value = 3680.0908
$ws.Range("M113").Value = -1510.0908

$ws.Range("H137").Value = 40700
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 441.1
$ws.Range("J12").Value = 629.7143
$ws.Range("L12").Value = 1889.1429
$ws.Range("N12").Value = -2235.1429

$ws.Range("H113").Value = 1796.0714
$ws.Range("I113").Value = 479.33334
$ws.Range("J113").Value = 2155.182
$ws.Range("K113").Value = 1438.00002
$ws.Range("L113").Value = 6465.545999999999
$ws.Range("M113").Value = 731.9999800000001
$ws.Range("N113").Value = -10805.546

$ws.Range("H131").Value = 19233380
$ws.Range("I131").Value = 50001170
$ws.Range("K131").Value = 150003510
$ws.Range("M131").Value = -149998470

$ws.Range("H132").Value = 19231632
$ws.Range("J132").Value = 1128.4445
$ws.Range("L132").Value = 10156.0005
$ws.Range("N132").Value = -15216.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6803.278
$ws.Range("I70").Value = 4206.5713
$ws.Range("J70").Value = 8455.727999999999
$ws.Range("K70").Value = 4206.5713
$ws.Range("L70").Value = 8455.727999999999
$ws.Range("M70").Value = -3936.5713
$ws.Range("N70").Value = -8995.727999999999

$ws.Range("H73").Value = 6803.278
$ws.Range("I73").Value = 4206.5713
$ws.Range("J73").Value = 8455.727999999999
$ws.Range("K73").Value = 4206.5713
$ws.Range("L73").Value = 8455.727999999999
$ws.Range("M73").Value = -3270.5713
$ws.Range("N73").Value = -10327.728

$ws.Range("H113").Value = 447037.22
$ws.Range("I113").Value = 668106.7
$ws.Range("K113").Value = 668106.7
$ws.Range("M113").Value = -665936.7

$ws.Range("H116").Value = 63745
$ws.Range("J116").Value = 63745
$ws.Range("L116").Value = 63745
$ws.Range("N116").Value = -72923

$ws.Range("H132").Value = 1660.2778
$ws.Range("I132").Value = 1356.3939
$ws.Range("K132").Value = 4069.1817
$ws.Range("M132").Value = -1539.1817

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5997
$ws.Range("I46").Value = 5997
$ws.Range("K46").Value = 5997
$ws.Range("M46").Value = -5809

$ws.Range("H61").Value = 144297.14
$ws.Range("I61").Value = 144297.14
$ws.Range("K61").Value = 144297.14
$ws.Range("M61").Value = -144095.14

$ws.Range("H68").Value = 3465.3333
$ws.Range("I68").Value = 2169.9285
$ws.Range("K68").Value = 2169.9285
$ws.Range("M68").Value = -1420.9285

$ws.Range("H71").Value = 3465.3333
$ws.Range("I71").Value = 2169.9285
$ws.Range("K71").Value = 10849.6425
$ws.Range("M71").Value = -7105.6425

$ws.Range("H99").Value = 22997.5
$ws.Range("I99").Value = 20663.334
$ws.Range("K99").Value = 20663.334
$ws.Range("M99").Value = -17668.334

$ws.Range("H113").Value = 144297.14
$ws.Range("I113").Value = 144297.14
$ws.Range("K113").Value = 144297.14
$ws.Range("M113").Value = -142127.14

$ws.Range("H132").Value = 9181.546
$ws.Range("I132").Value = 9231.561
$ws.Range("J132").Value = 9035.071
$ws.Range("K132").Value = 27694.683
$ws.Range("L132").Value = 27105.213
$ws.Range("M132").Value = -25164.683
$ws.Range("N132").Value = -32165.213

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H113").Value = 1357
$ws.Range("I113").Value = 992.913
$ws.Range("J113").Value = 3450.5
$ws.Range("K113").Value = 2978.739
$ws.Range("L113").Value = 10351.5
$ws.Range("M113").Value = -808.739
$ws.Range("N113").Value = -14691.5

$ws.Range("H122").Value = 4491.879
$ws.Range("I122").Value = 3677.9565
$ws.Range("K122").Value = 11033.8695
$ws.Range("M122").Value = -8583.869499999999

$ws.Range("H136").Value = 7789.8335
$ws.Range("J136").Value = 8721.111000000001
$ws.Range("L136").Value = 26163.333
$ws.Range("N136").Value = -31263.333
